# DEV-2004 Simplify Import Tasks - Remove Logstash
# Update XLSX sample (order_import_sample.xlsx):
#  - rename the sample SKUs used on rows 2-4
#  - bump the sample order_ref on rows 2-4
#  - leave the sheet's cursor on K9 (matches the saved selection in the sample)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# sku column (L) - rename sample SKUs: test-sku-1/2/3 -> test-sku-6/7/8
$ws.Range("L2").Value = "test-sku-6"
$ws.Range("L3").Value = "test-sku-7"
$ws.Range("L4").Value = "test-sku-8"

# order_ref column (A) - sample order reference number updated on every row
$ws.Range("A2").Value = 323456
$ws.Range("A3").Value = 323456
$ws.Range("A4").Value = 323456

# Move/save the active selection to K9, as recorded in the updated sample
$ws.Range("K9").Select()
